$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2-5: date serial 45204 -> 45207 (2023-10-05 -> 2023-10-08)
$ws.Range("C2:C5").Value = 45207
